$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.229.32'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '2.366.89'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.32'
$ws.Range("E5").Value = '  +1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.80'
$ws.Range("E6").Value = '  +2.69%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("D9").Value = '2.363.99'
$ws.Range("E9").Value = '  +1.24%  '

$ws.Range("E10").Value = '  +3.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  +1.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.351'
$ws.Range("E13").Value = '  +3.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.60'
$ws.Range("E14").Value = '  +4.82%  '

$ws.Range("E15").Value = '  +7.91%  '

$ws.Range("D16").Value = '2.789.77'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").Value = '61.301.25'

$ws.Range("D18").Value = '2.361.94'
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.03'
$ws.Range("E19").Value = '  +5.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.15'
$ws.Range("E20").Value = '  +2.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.33'
$ws.Range("E21").Value = '  +1.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.63'
$ws.Range("E22").Value = '  +1.45%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.26'
$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("E25").Value = '  -6.95%  '

$ws.Range("E26").Value = '  +3.23%  '

$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '536.90'
$ws.Range("E28").Value = '  +7.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.24'
$ws.Range("E30").Value = '  +4.09%  '

$ws.Range("D31").Value = '0.0₃0907'
$ws.Range("E31").Value = '  +2.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  +0.76%  '

$ws.Range("E33").Value = '  +2.05%  '

$ws.Range("E34").Value = '  +3.37%  '

$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.63'
$ws.Range("E37").Value = '  +8.10%  '

$ws.Range("E38").Value = '  +2.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("E39").Value = '  +6.31%  '

$ws.Range("E40").Value = '  +2.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.49'
$ws.Range("E41").Value = '  +1.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '145.86'
$ws.Range("E42").Value = '  +6.12%  '

$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.47'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.37'
$ws.Range("E45").Value = '  +4.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  +6.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.60'
$ws.Range("E47").Value = '  +2.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0529'
$ws.Range("E48").Value = '  +3.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.11'
$ws.Range("E49").Value = '  +4.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.580'
$ws.Range("E50").Value = '  +2.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0903'
$ws.Range("E51").Value = '  +0.77%  '
